$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.726
$ws.Range("B3").Value = 6.201
$ws.Range("C5").Value = -12.836
$ws.Range("E5").Value = 13.034
$ws.Range("E9").Value = 13.152
$ws.Range("E11").Value = 13.071
$ws.Range("B14").Value = 6.449
$ws.Range("B21").Value = 6.6
$ws.Range("E21").Value = 12.694
$ws.Range("B23").Value = 6.610000000000001
$ws.Range("B25").Value = 5.986
